# Add a dummy "Al" (Aluminum) data row to the "properties" sheet.
# The new row (35) duplicates all of row 7's (Gallium, "Ga") numeric
# property values/formatting verbatim - only the Symbol cell changes,
# from "Ga" to "Al".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number formats, fills, fonts, etc.) from the Gallium
# row (7) down onto the new row (35) so the new row's cells line up
# with the existing per-column style reused across the table.
$ws.Range("A7:AB7").Copy()
$ws.Range("A35").PasteSpecial(-4122)

# Copy the underlying values too (same data used for the dummy row).
for ($col = 2; $col -le 28; $col++) {
    $src = $ws.Cells.Item(7, $col)
    $dst = $ws.Cells.Item(35, $col)
    $dst.Value = $src.Value()
}

# Only the element symbol differs for this dummy row.
$ws.Cells.Item(35, 1).Value = "Al"

# Leave the sheet selection on the newly added row, like the author did.
[void]$ws.Range("B35:AB35").Select()
